$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- page-rank sheet: add two new blocks of columns (K:M and O:Q) ---

# Row 1 headers
$ws1.Range("K1").Value = "map - split"
$ws1.Range("O1").Value = "mapPar-split"

# Row 2 sub-headers
$ws1.Range("K2").Value = "computation"
$ws1.Range("L2").Value = "total"
$ws1.Range("M2").Value = "merge"

$ws1.Range("O2").Value = "computation"
$ws1.Range("P2").Value = "split"
$ws1.Range("Q2").Value = "merge"

# Row 3 data
$ws1.Range("K3").Value = 0.31519290773303998
$ws1.Range("L3").Value = 1.0316822960826699
$ws1.Range("M3").Value = 0.073807045950577493

$ws1.Range("O3").Value = 0.30653552513797999
$ws1.Range("P3").Value = 1.01931440411138
$ws1.Range("Q3").Value = 0.059332231192476902

# Row 4 data
$ws1.Range("K4").Value = 0.40530568339822298
$ws1.Range("L4").Value = 1.10064360867004
$ws1.Range("M4").Value = 0.0648689293637362

$ws1.Range("O4").Value = 0.41344602392760998
$ws1.Range("P4").Value = 1.14696613516254
$ws1.Range("Q4").Value = 0.066403041379433897

# Row 5 data
$ws1.Range("K5").Value = 0.41530198753994702
$ws1.Range("L5").Value = 1.22933320478174
$ws1.Range("M5").Value = 0.13133624653926801

$ws1.Range("O5").Value = 0.54558571791810595
$ws1.Range("P5").Value = 1.2081992285014
$ws1.Range("Q5").Value = 0.059687638535608503

# Row 6 data
$ws1.Range("K6").Value = 0.40404048997061198
$ws1.Range("L6").Value = 1.12248556393976
$ws1.Range("M6").Value = 0.065359977541281297

$ws1.Range("O6").Value = 0.36761501822880399
$ws1.Range("P6").Value = 1.0544482159213
$ws1.Range("Q6").Value = 0.058096812582176402

# Row 7 data
$ws1.Range("K7").Value = 0.38736233873155501
$ws1.Range("L7").Value = 1.22485658462995
$ws1.Range("M7").Value = 0.057353859997862799

$ws1.Range("O7").Value = 0.41728839320907302
$ws1.Range("P7").Value = 1.3216601870992799
$ws1.Range("Q7").Value = 0.066089696873454104

# Row 9 averages
$ws1.Range("K9").Formula = "=AVERAGE(K3:K7)"
$ws1.Range("O9").Formula = "=AVERAGE(O3:O7)"

# --- update selections: spmv sheet first so page-rank ends up as the active tab ---
$ws2.Range("E8").Select()
$ws1.Activate()
$ws1.Range("L13").Select()
